$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column B holds plain-text dates ("2024-11-09"). Assigning a date-shaped
# string straight to .Value makes Excel re-interpret it as a real date
# (serial number + a new date number-format style), which would change
# the cell's style index. To keep the values as literal text with their
# original style untouched, stash a copy of the existing (text) format in
# a scratch cell, then re-apply it via PasteSpecial(Formats) after setting
# the new literal text value.

$holder = $ws.Cells.Item(1, 10)
$template = $ws.Cells.Item(2, 2)
$template.Copy()
$holder.PasteSpecial(-4122)  # xlPasteFormats

for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = "2024-11-18"
    $holder.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats - restore original style
}

$holder.Clear()
$excel.CutCopyMode = $false
